$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-7 from 45204 to 45207 (date serials)
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45207
}
